# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 6080
$sheet1.Range("F6").Value = 51
$sheet1.Range("F9").Value = 52
$sheet1.Range("F10").Value = 64
$sheet1.Range("F14").Value = 608
$sheet1.Range("F15").Value = 3078
$sheet1.Range("F17").Value = 177
$sheet1.Range("F18").Value = 1727
$sheet1.Range("F19").Value = 18

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 6080
$sheet4.Range("F6").Value = 51
$sheet4.Range("F10").Value = 52
$sheet4.Range("F11").Value = 64
$sheet4.Range("F15").Value = 608
$sheet4.Range("F16").Value = 3078
$sheet4.Range("F18").Value = 177
$sheet4.Range("F19").Value = 1727
$sheet4.Range("F20").Value = 18
